$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 58; existing rows 58-64 shift down to 59-65
$ws.Rows.Item(58).Insert()

# Populate the new row 58 with this week's record (copy of the constant
# columns used throughout this sheet, plus the new weekly values)
$ws.Cells.Item(58, 1).Value = 7
$ws.Cells.Item(58, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(58, 3).Value = "Ñuble"
$ws.Cells.Item(58, 4).Value = 45223
$ws.Cells.Item(58, 5).Value = 16
$ws.Cells.Item(58, 6).Value = 300000000
$ws.Cells.Item(58, 7).Value = "Espárragos"
$ws.Cells.Item(58, 8).Value = "Sin especificar"
$ws.Cells.Item(58, 9).Value = "Primera"
$ws.Cells.Item(58, 10).Value = 1000
$ws.Cells.Item(58, 11).Value = 1000
$ws.Cells.Item(58, 12).Value = 1200
$ws.Cells.Item(58, 13).Value = 1100
$ws.Cells.Item(58, 14).Value = "$/kilo"
$ws.Cells.Item(58, 15).Value = "Región de Ñuble"
$ws.Cells.Item(58, 16).Value = 1100
$ws.Cells.Item(58, 17).Value = 1
$ws.Cells.Item(58, 18).Value = "Hortaliza"
